# PyPanel_class-subclass_listing.xlsx update:
#  - trim the "turtle" method list from the full CircuitPython `turtle`
#    module reference down to the subset PyPanel actually subclasses
#    (rows 4:44 in column K), drop the now-unused tail rows (45:91 lost
#    their column-K entries, and rows 53:91 are gone entirely)
#  - rename the "panel (via PyBadger)" column/header to "panel"
#  - rename the "joystick (PyBadge only)" cell to "joystick (PyGamer only)"
#  - re-center/re-scale the print setup

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- shrink the sheet: drop rows 53:91 outright, then blank out the
#     column-K tail that used to run through row 91 but now stops at 44 ---
$ws.Rows("53:91").Delete()
$ws.Range("K45:K52").ClearContents()

# --- column K (turtle-derived method list), rows 4:44, replaced with the
#     curated PyPanel subset, alphabetised ---
$turtleMethods = @(
  "addshape","back","backward","bk","circle","clear","degrees","dot",
  "down","fd","forward","getpen","goto","heading","home","ht","isdown",
  "left","lt","pd","pencolor","pendown","pensize","penup","position",
  "pu","radians","right","rt","seth","setheading","setpos","setposition",
  "setx","sety","st","turtlesize","up","width","xcor","ycor"
)
for ($i = 0; $i -lt $turtleMethods.Length; $i++) {
  $ws.Cells.Item(4 + $i, 11).Value = $turtleMethods[$i]
}

# --- header renames ---
$ws.Range("E3").Value = "panel"
$ws.Range("E17").Value = "joystick (PyGamer only)"

# --- print setup: center on page both ways, rescale to fit ---
$ws.PageSetup.CenterHorizontally = $true
$ws.PageSetup.CenterVertically = $true
$ws.PageSetup.Zoom = 69

# --- restore the view / selection state ---
$ws.Activate()
$ws.Range("E29").Select()
